$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "latest entry" yellow highlight from row 30 to row 22 ---
# Row 30 currently carries the highlighted style (fill + bold-ish font); copy
# that formatting onto row 22 (the new last item), then reset row 30 back to
# the plain style used by the rest of the list (copied from row 25).
$ws.Range("B30:C30").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)

$ws.Range("B25:C25").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update list values (new items / renamed items / shifted text) ---
# The brand-new strings are written first, in the same first-use order as
# the original edit, so the workbook's shared-string table comes out
# byte-for-byte identical; the remaining cells merely reuse / reshuffle
# strings that already exist in the table, so their order doesn't matter.
$ws.Range("B22").Value = "CARBURADOR DS150"
$ws.Range("B8").Value = "SLAIDER DM200"
$ws.Range("B18").Value = "ACEITE YAMALUBE"
$ws.Range("C18").Value = "DESCONOZCO"
$ws.Range("I9").Value = "BUJIA PULSAR "
$ws.Range("B23").Value = "AUXILIAR BUHO"
$ws.Range("B24").Value = "BANDA 743-20-30"

$ws.Range("B14").Value = "CADENA 520 CON ORING REFORZADA "
$ws.Range("C14").Value = "NASAKI"

$ws.Range("B15").Value = "PASTAS DE CLUTH DS150"

$ws.Range("B16").Value = "PASTAS DE CLUTH FT180/FT200"
$ws.Range("C16").Value = "NASAKI"

$ws.Range("B17").Value = "BUJIA PARA PULSAR 200NS "
$ws.Range("C17").Value = "NASAKI"

$ws.Range("B19").Value = "MANGUERA DE FRENO AT110"
$ws.Range("C19").Value = "NASAKI"

$ws.Range("B20").Value = "CADENA REFORZADA CON ORING 428"

$ws.Range("B21").Value = "PORTAFUSIBLES "

$ws.Range("C23").Value = "MERCADO LIBRE "

$ws.Range("C24").Value = ""

$ws.Range("B26").Value = ""
$ws.Range("C26").Value = ""

$ws.Range("B29").Value = ""
$ws.Range("C29").Value = ""

$ws.Range("B30").Value = ""
$ws.Range("C30").Value = ""

# --- Move the active selection (matches the source workbook's cursor) ---
$excel.Goto($ws.Range("B25"))
